$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.476.75'
$ws.Range("E2").Value = '  -2.21%  '
$ws.Range("D3").Value = '2.521.86'
$ws.Range("E3").Value = '  -5.78%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = "'575.39"
$ws.Range("E5").Value = '  -4.05%  '
$ws.Range("D6").Value = "'169.28"
$ws.Range("E6").Value = '  -3.41%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -3.40%  '
$ws.Range("D9").Value = '2.517.90'
$ws.Range("E9").Value = '  -5.94%  '
$ws.Range("E10").Value = '  -3.26%  '
$ws.Range("E11").Value = '  -0.72%  '
$ws.Range("E12").Value = '  -3.91%  '
$ws.Range("D13").Value = "'4.79"
$ws.Range("E13").Value = '  -4.09%  '
$ws.Range("D14").Value = '2.979.48'
$ws.Range("E14").Value = '  -6.07%  '
$ws.Range("D15").Value = '70.265.11'
$ws.Range("E15").Value = '  -2.39%  '
$ws.Range("D16").Value = "'0.0000179"
$ws.Range("E16").Value = '  -3.02%  '
$ws.Range("D17").Value = "'24.86"
$ws.Range("E17").Value = '  -5.43%  '
$ws.Range("D18").Value = '2.516.53'
$ws.Range("E18").Value = '  -6.06%  '
$ws.Range("D19").Value = "'11.48"
$ws.Range("E19").Value = '  -6.34%  '
$ws.Range("D20").Value = "'7.53"
$ws.Range("E20").Value = '  -8.22%  '
$ws.Range("D21").Value = "'355.85"
$ws.Range("E21").Value = '  -4.65%  '
$ws.Range("E22").Value = '  -6.24%  '
$ws.Range("E23").Value = '  -4.65%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").Value = "'69.08"
$ws.Range("E25").Value = '  -4.17%  '
$ws.Range("D26").Value = "'4.05"
$ws.Range("E26").Value = '  -7.08%  '
$ws.Range("D27").Value = "'9.18"
$ws.Range("E27").Value = '  -6.39%  '
$ws.Range("D28").Value = '2.648.52'
$ws.Range("E28").Value = '  -5.97%  '
$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("D30").Value = '0.0₃0907'
$ws.Range("E30").Value = '  -6.77%  '
$ws.Range("E31").Value = '  -3.34%  '
$ws.Range("D32").Value = "'478.07"
$ws.Range("E32").Value = '  -4.76%  '
$ws.Range("E33").Value = '  -3.71%  '
$ws.Range("E34").Value = '  -3.96%  '
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = '  -0.23%  '
$ws.Range("D36").Value = "'158.45"
$ws.Range("E36").Value = '  -2.59%  '
$ws.Range("E37").Value = '  +3.99%  '
$ws.Range("D38").Value = "'18.83"
$ws.Range("D39").Value = "'18.53"
$ws.Range("E39").Value = '  -5.37%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("E41").Value = '  -5.89%  '
$ws.Range("E42").Value = '  -7.59%  '
$ws.Range("D43").Value = "'0.318"
$ws.Range("E43").Value = '  -4.75%  '
$ws.Range("D44").Value = "'4.70"
$ws.Range("E44").Value = '  -6.28%  '
$ws.Range("E45").Value = '  -6.68%  '
$ws.Range("E46").Value = '  -3.10%  '
$ws.Range("D47").Value = "'142.44"
$ws.Range("E47").Value = '  -8.98%  '
$ws.Range("D48").Value = "'3.52"
$ws.Range("E48").Value = '  -5.93%  '
$ws.Range("E49").Value = '  -7.14%  '
$ws.Range("E50").Value = '  -7.65%  '
$ws.Range("D51").Value = "'0.593"
$ws.Range("E51").Value = '  -2.11%  '
